# Add a new weekly price record for "Hortaliza, Terminal Hortofrutícola Agro
# Chillán - Coliflor" by inserting a new row above the current row 512. This
# shifts the existing rows 512-557 down to 513-558 (preserving all of their
# data), and the newly inserted row 512 is populated with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 512; Excel shifts rows 512:557 down to 513:558.
$ws.Rows.Item(512).Insert()

# Populate the new row 512 with the new market record.
$ws.Cells.Item(512, 1).Value  = 7
$ws.Cells.Item(512, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(512, 3).Value  = "Ñuble"
$ws.Cells.Item(512, 4).Value  = 45166
$ws.Cells.Item(512, 5).Value  = 16
$ws.Cells.Item(512, 6).Value  = 100112008
$ws.Cells.Item(512, 7).Value  = "Coliflor"
$ws.Cells.Item(512, 8).Value  = "Sin especificar"
$ws.Cells.Item(512, 9).Value  = "Primera"
$ws.Cells.Item(512, 10).Value = 180
$ws.Cells.Item(512, 11).Value = 1000
$ws.Cells.Item(512, 12).Value = 1000
$ws.Cells.Item(512, 13).Value = 1000
$ws.Cells.Item(512, 14).Value = "`$/unidad"
$ws.Cells.Item(512, 15).Value = "Región del Maule"
$ws.Cells.Item(512, 16).Value = 1000
$ws.Cells.Item(512, 17).Value = 1
$ws.Cells.Item(512, 18).Value = "Hortaliza"
